$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.558.19'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '2.617.09'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''593.76'
$ws.Range('E5').Value = '  -1.65%  '
$ws.Range('D6').Value = '''150.21'
$ws.Range('E6').Value = '  +2.69%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '''0.587'
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('D10').Value = '''5.69'
$ws.Range('E10').Value = '  +2.05%  '
$ws.Range('E11').Value = '  +3.07%  '
$ws.Range('E12').Value = '  -1.15%  '
$ws.Range('D13').Value = '''27.60'
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('D14').Value = '3.089.23'
$ws.Range('E14').Value = '  -0.88%  '
$ws.Range('D15').Value = '63.410.39'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('D16').Value = '''0.0000151'
$ws.Range('E16').Value = '  +3.01%  '
$ws.Range('D17').Value = '2.629.05'
$ws.Range('E17').Value = '  -0.90%  '
$ws.Range('D18').Value = '''12.32'
$ws.Range('E18').Value = '  +7.49%  '
$ws.Range('D19').Value = '''4.65'
$ws.Range('E19').Value = '  +1.74%  '
$ws.Range('D20').Value = '''346.10'
$ws.Range('E20').Value = '  +0.82%  '
$ws.Range('D21').Value = '''6.85'
$ws.Range('E21').Value = '  -1.24%  '
$ws.Range('D22').Value = '''0.996'
$ws.Range('E22').Value = '  -0.42%  '
$ws.Range('D23').Value = '''5.72'
$ws.Range('E23').Value = '  +2.72%  '
$ws.Range('D24').Value = '''66.30'
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('D25').Value = '''1.72'
$ws.Range('E25').Value = '  +10.54%  '
$ws.Range('D26').Value = '''9.20'
$ws.Range('E26').Value = '  +1.25%  '
$ws.Range('D27').Value = '''1.67'
$ws.Range('E27').Value = '  -1.40%  '
$ws.Range('D28').Value = '''562.30'
$ws.Range('E28').Value = '  -2.71%  '
$ws.Range('D29').Value = '''8.21'
$ws.Range('E29').Value = '  +3.66%  '
$ws.Range('D30').Value = '''0.162'
$ws.Range('E30').Value = '  -0.68%  '
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('E32').Value = '  -0.57%  '
$ws.Range('D33').Value = '0.0₃0843'
$ws.Range('E33').Value = '  +2.37%  '
$ws.Range('D34').Value = '''1.75'
$ws.Range('E34').Value = '  -0.50%  '
$ws.Range('D35').Value = '''5.23'
$ws.Range('E35').Value = '  +0.45%  '
$ws.Range('D36').Value = '''168.44'
$ws.Range('E36').Value = '  +0.89%  '
$ws.Range('D37').Value = '''0.408'
$ws.Range('E37').Value = '  +0.52%  '
$ws.Range('D38').Value = '''1.00'
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('D40').Value = '''19.32'
$ws.Range('E40').Value = '  +1.00%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('D42').Value = '''166.94'
$ws.Range('E42').Value = '  -1.14%  '
$ws.Range('D43').Value = '''39.93'
$ws.Range('E43').Value = '  -0.16%  '
$ws.Range('D44').Value = '''3.91'
$ws.Range('E44').Value = '  +3.77%  '
$ws.Range('D45').Value = '''0.0596'
$ws.Range('E45').Value = '  +4.72%  '
$ws.Range('D46').Value = '''21.39'
$ws.Range('E46').Value = '  -3.38%  '
$ws.Range('D47').Value = '''0.626'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('E48').Value = '  +1.39%  '
$ws.Range('E49').Value = '  +4.40%  '
$ws.Range('D50').Value = '''0.0962'
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('D51').Value = '''19.36'
$ws.Range('E51').Value = '  +3.18%  '
